$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(1)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# The target paragraph is the 3rd paragraph of the content placeholder:
# "El ranking més alt al que hem arribat dins de Bronze ha estat aproximadament la posició 1150."
$para = $tr.Paragraphs(3)

# 1) Change the number "1150" -> "750" (this splits the single run into
#    "...la posició " / "750" / "." - mirroring how PowerPoint splits a run
#    when only part of its text is edited).
$numStart = $para.Text.IndexOf("1150") + 1
$numRange = $para.Characters($numStart, 4)
$numRange.Text = "750"

# 2) Split "posició " into its own run, leaving the leading
#    "El ranking ... aproximadament la " text as the first run.
$posStart = $para.Text.IndexOf("posició ") + 1
$posRange = $para.Characters($posStart, 8)
$posRange.Text = "posició "

# 3) Split the trailing "." into its own run.
$dotRange = $para.Characters($para.Length, 1)
$dotRange.Text = "."
